# Update the build timestamp embedded in the workbook from
# "January 30 2026 16.19.47 EST" to "February 02 2026 12.49.33 EST"
# across the "About" sheet (version/citation text) and the
# "Boundaries and methane sources" sheet (build_version column).

$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

$oldVersion = "mines - January 30 (built on $oldStamp)"
$newVersion = "mines - January 30 (built on $newStamp)"

# --- "About" sheet ---
$wsAbout = $wb.Worksheets.Item("About")

$a2 = $wsAbout.Range("A2")
$a2.Value = $a2.Value().Replace($oldStamp, $newStamp)

$a6 = $wsAbout.Range("A6")
$a6.Value = $a6.Value().Replace($oldStamp, $newStamp)

# --- "Boundaries and methane sources" sheet ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 13; $row++) {
    $cell = $wsData.Cells.Item($row, 19)  # column S = build_version
    if ($cell.Value() -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
